$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("input_concentrations")

# Header renamed from "ser_num" to "series"
$ws.Range("C2").Value = "series"

# Restore this sheet as active and update the selected cell to H11
$ws.Activate()
$ws.Range("H11").Select()
